# Fixed sensor data extraction
# Re-order the attribute/type rows in the main event attribute info table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("process_model_id", "str"),
    @("response_status_code", "float"),
    @("SubProcessID", "str"),
    @("case", "str"),
    @("complete_service_time", "str"),
    @("parameters", "dict"),
    @("lifecycle:transition", "str"),
    @("human_workstation_green_button_pressed", "float"),
    @("org:resource", "str"),
    @("requested_service_url", "str"),
    @("lifecycle:state", "str"),
    @("identifier:id", "str"),
    @("current_task", "str"),
    @("event_id", "str"),
    @("unsatisfied_condition_description", "str"),
    @("case:concept:name", "str"),
    @("concept:name", "str"),
    @("planned_operation_time", "str"),
    @("operation_end_time", "datetime"),
    @("time:timestamp", "datetime")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
